$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 142.9073533333333
$ws.Range("H2").Value = 428.72206
$ws.Range("I2").Value = 0.5576664151504187
$ws.Range("J2").Value = 0.5576664151504188
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.8369776666666665
$ws.Range("N2").Value = 2.510933
$ws.Range("O2").Value = 0.0694586718035551
$ws.Range("P2").Value = 0.06945867180355511
$ws.Range("Q2").Value = 119.6102631424422
$ws.Range("R2").Value = 1076.49236828198
$ws.Range("S2").Value = 0.03873476850579804
$ws.Range("T2").Value = 0.03873476850579805
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 142.9073533333333
$ws.Range("H3").Value = 428.72206
$ws.Range("I3").Value = 0.5576664151504187
$ws.Range("J3").Value = 0.5576664151504188
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 2.427350333333333
$ws.Range("N3").Value = 7.282051
$ws.Range("O3").Value = 0.2014397000898671
$ws.Range("P3").Value = 0.2014397000898671
$ws.Range("Q3").Value = 346.8862117494511
$ws.Range("R3").Value = 3121.97590574506
$ws.Range("S3").Value = 0.1123361554180917
$ws.Range("T3").Value = 0.1123361554180917
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 142.9073533333333
$ws.Range("H4").Value = 428.72206
$ws.Range("I4").Value = 0.5576664151504187
$ws.Range("J4").Value = 0.5576664151504188
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 8.785681666666667
$ws.Range("N4").Value = 26.357045
$ws.Range("O4").Value = 0.7291016281065776
$ws.Range("P4").Value = 0.7291016281065776
$ws.Range("Q4").Value = 1255.538514212522
$ws.Range("R4").Value = 11299.8466279127
$ws.Range("S4").Value = 0.4065954912265289
$ws.Range("T4").Value = 0.406595491226529
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 63.967809
$ws.Range("H5").Value = 191.903427
$ws.Range("I5").Value = 0.2496211559306514
$ws.Range("J5").Value = 0.2496211559306514
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 0.8369776666666665
$ws.Range("N5").Value = 2.510933
$ws.Range("O5").Value = 0.0694586718035551
$ws.Range("P5").Value = 0.06945867180355511
$ws.Range("Q5").Value = 53.53962751859898
$ws.Range("R5").Value = 481.8566476673909
$ws.Range("S5").Value = 0.01733835394501116
$ws.Range("T5").Value = 0.01733835394501117
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 63.967809
$ws.Range("H6").Value = 191.903427
$ws.Range("I6").Value = 0.2496211559306514
$ws.Range("J6").Value = 0.2496211559306514
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 2.427350333333333
$ws.Range("N6").Value = 7.282051
$ws.Range("O6").Value = 0.2014397000898671
$ws.Range("P6").Value = 0.2014397000898671
$ws.Range("Q6").Value = 155.272282498753
$ws.Range("R6").Value = 1397.450542488777
$ws.Range("S6").Value = 0.05028361078675637
$ws.Range("T6").Value = 0.05028361078675638
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 63.967809
$ws.Range("H7").Value = 191.903427
$ws.Range("I7").Value = 0.2496211559306514
$ws.Range("J7").Value = 0.2496211559306514
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 8.785681666666667
$ws.Range("N7").Value = 26.357045
$ws.Range("O7").Value = 0.7291016281065776
$ws.Range("P7").Value = 0.7291016281065776
$ws.Range("Q7").Value = 562.000806788135
$ws.Range("R7").Value = 5058.007261093215
$ws.Range("S7").Value = 0.1819991911988838
$ws.Range("T7").Value = 0.1819991911988838
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 49.38440333333333
$ws.Range("H8").Value = 148.15321
$ws.Range("I8").Value = 0.1927124289189298
$ws.Range("J8").Value = 0.1927124289189298
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.8369776666666665
$ws.Range("N8").Value = 2.510933
$ws.Range("O8").Value = 0.0694586718035551
$ws.Range("P8").Value = 0.06945867180355511
$ws.Range("Q8").Value = 41.33364267165888
$ws.Range("R8").Value = 372.00278404493
$ws.Range("S8").Value = 0.01338554935274589
$ws.Range("T8").Value = 0.01338554935274589
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 49.38440333333333
$ws.Range("H9").Value = 148.15321
$ws.Range("I9").Value = 0.1927124289189298
$ws.Range("J9").Value = 0.1927124289189298
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.427350333333333
$ws.Range("N9").Value = 7.282051
$ws.Range("O9").Value = 0.2014397000898671
$ws.Range("P9").Value = 0.2014397000898671
$ws.Range("Q9").Value = 119.8732478926345
$ws.Range("R9").Value = 1078.85923103371
$ws.Range("S9").Value = 0.03881993388501906
$ws.Range("T9").Value = 0.03881993388501907
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 49.38440333333333
$ws.Range("H10").Value = 148.15321
$ws.Range("I10").Value = 0.1927124289189298
$ws.Range("J10").Value = 0.1927124289189298
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.785681666666667
$ws.Range("N10").Value = 26.357045
$ws.Range("O10").Value = 0.7291016281065776
$ws.Range("P10").Value = 0.7291016281065776
$ws.Range("Q10").Value = 433.8756469849389
$ws.Range("R10").Value = 3904.88082286445
$ws.Range("S10").Value = 0.1405069456811648
$ws.Range("T10").Value = 0.1405069456811649
